$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell values - headers + the 4 credential rows (new BEES lines added)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = ""
$ws.Range("B1").Value = "NOMBRE DE USUARIO"
$ws.Range("C1").Value = "CONTRASEÑA ENCRIPTADA"

$ws.Range("A2").Value = "MAXICONSUMO"
$ws.Range("B2").Value = "orlando.piccinini@gmail.com"
$ws.Range("C2").Value = "26TBDyfahsU=*kVpqWxsj5NswRl9AYd9qVw==*Ok3UyLr0xbq0hjgxP0lomw==*NepodkMzbVYH+ModcC9Sgw=="

$ws.Range("A3").Value = "LA SERENÍSIMA"
$ws.Range("B3").Value = "orlando.piccinini@gmail.com"
$ws.Range("C3").Value = "7bwYIZivWO4dVRM=*mWprFGvNLOSV5+crKLSdxQ==*lhWQVj7MPOLpUA4Q8ciZ3Q==*q7T1Qa4n1CiByRrvIilD0g=="

$ws.Range("A4").Value = "BEES (GRAL. ALVEAR)"
$ws.Range("B4").Value = "2625404916"
$ws.Range("C4").Value = "ddALT7mE3UQwH+DgsjKs900=*XX5AKld2q0bF52DUbAnM8w==*75+mUm0SnqACoKlKdheVgQ==*SsCetcMJvyJU29lez1AvQA=="

$ws.Range("A5").Value = "BEES (SAN RAFAEL)"
$ws.Range("B5").Value = "1158108611"
$ws.Range("C5").Value = "BXDa7+m0Z3fg*sOGWT/rynsrO5dtvUe7tgw==*qAs8nFu3/+r9znWBh8KsUA==*eescdr4iBwPxuRun34hgZg=="

# ---------------------------------------------------------------------------
# 2) Borders - every cell in A1:C5 gets a thin box border. Order matters: we
#    group cells by their CURRENT shared style before touching them so the
#    engine doesn't strand/duplicate style records.
# ---------------------------------------------------------------------------

# B3 carried a leftover (no-border, underlined) placeholder style - normalize
# it before bordering so it merges with the rest of the plain cells.
$ws.Range("B3").Style = "Normal"

# B1 & B2 both still carry the ORIGINAL "Hipervínculo" style at this point -
# border them together first so that style moves as a whole (no orphan).
$ws.Range("B1:B2").Borders.LineStyle = 1

# B1 is actually a plain header cell (no hyperlink) in the final layout.
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Borders.LineStyle = 1

# Border the remaining plain cells.
$ws.Range("A1:A5").Borders.LineStyle = 1
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C3:C5").Borders.LineStyle = 1
$ws.Range("B3:B5").Borders.LineStyle = 1

# C2 must look like B2 (hyperlink style) plus the border.
$ws.Range("C2").Style = "Hipervínculo"
$ws.Range("C2").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3) Hyperlinks - B2 and C2 both link to the same mailto address.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:orlando.piccinini@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:orlando.piccinini@gmail.com")

# ---------------------------------------------------------------------------
# 4) Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 111.5703125
$ws.Columns.Item(6).ColumnWidth = 113.7109375
$ws.Columns.Item(7).ColumnWidth = 31.42578125
$ws.Columns.Item(8).ColumnWidth = 101.5703125

# ---------------------------------------------------------------------------
# 5) Selection
# ---------------------------------------------------------------------------
$ws.Range("B8").Select()

Write-Host "done"
